# Windows Compatible version in arddata change to the port connected to arduino
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unprotect the workbook (drops the stray <workbookProtection/> element)
$wb.Unprotect()

# The port connected to the arduino changed from "TD" to "TW" for the
# first batch of entries (rows 1-4). Rows 5-6 are left untouched.
$ws.Range("C1").Value = "43 TR IO TW"
$ws.Range("C2").Value = "43 TR IO TW"
$ws.Range("C3").Value = "43 TR IO TW"
$ws.Range("C4").Value = "43 TR IO TW"

# Append new rows of arduino/device data (rows 7-12).
# Column B holds numeric-looking IDs that must stay text, so they are
# entered with a leading apostrophe (same as typing them in Excel's UI).
$ws.Range("A7").Value = "LOVE"
$ws.Range("B7").Formula = "'219311046"
$ws.Range("C7").Value = "43 TR IO TD"

$ws.Range("A8").Value = "LOVE"
$ws.Range("B8").Formula = "'219311046"
$ws.Range("C8").Value = "43 TR IO TD"

$ws.Range("A9").Value = "LOVE"
$ws.Range("B9").Formula = "'219311046"
$ws.Range("C9").Value = "43 TR IO TD"

$ws.Range("A10").Value = "LOVE"
$ws.Range("B10").Formula = "'219311046"
$ws.Range("C10").Value = "43 TR IO TD"

$ws.Range("A11").Value = "Love Lakhwani"
$ws.Range("B11").Formula = "'219311046"
$ws.Range("C11").Value = "1D 72 B0 04"

$ws.Range("A12").Value = "Utkarsh Triphati"
$ws.Range("B12").Formula = "'211015048"
$ws.Range("C12").Value = "AD B2 D3 04"

# Leave the selection on the last-touched cell, matching the saved view state
$ws.Range("F10").Select()
